# The post at row 581 ("「人間は一人ひとり、特別な才がある。だから自分を他人と比べるのは
# やめよう」") was removed from the source data. Deleting its entire row shifts every
# subsequent row up by one, matching the rest of the diff (rows 582-675 -> 581-674) and
# the updated used-range dimension (A1:C675 -> A1:C674).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(581).Delete()
